# Actualización automática del mapa: agrega el nuevo registro (fila 85)
# al final de la hoja "PEBCOM", replicando el formato de las filas previas
# (columnas de texto/fecha/identificador almacenadas como texto,
# columnas numéricas I/M/N como números).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85

function Set-TextCell($r, $c, $val) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $row 1  "-566"
Set-TextCell $row 2  "8/25/2025"
Set-TextCell $row 3  "Asuncion 2710"
Set-TextCell $row 4  "15"
Set-TextCell $row 5  "809171103"
Set-TextCell $row 6  "PEBCOM"
Set-TextCell $row 7  "Pendiente"
Set-TextCell $row 8  "Cambio"
$ws.Cells.Item($row, 9).Value = 0
Set-TextCell $row 10 "Cambio"
Set-TextCell $row 11 "Sin equipos"
Set-TextCell $row 12 "Terminal"
$ws.Cells.Item($row, 13).Value = -58.494789
$ws.Cells.Item($row, 14).Value = -34.59082
Set-TextCell $row 15 "Paternal"
Set-TextCell $row 16 "Capital Norte"
